{"js": "// Update the date line and the 24 multiplication expressions in the table,\n// per the authoring diff. Every old string below is unique within the\n// document, so a matchCase/matchWholeWord search finds exactly one hit\n// each; we then replace that hit's text in place so the run's existing\n// formatting (fonts/size) is preserved.\nconst replacements = [\n  [\"2025-10-10 Friday\", \"2025-10-11 Saturday\"],\n  [\"887\u00d78=\", \"201\u00d76=\"],\n  [\"953\u00d75=\", \"257\u00d76=\"],\n  [\"463\u00d74=\", \"411\u00d77=\"],\n  [\"141\u00d74=\", \"191\u00d72=\"],\n  [\"579\u00d75=\", \"580\u00d72=\"],\n  [\"169\u00d75=\", \"643\u00d76=\"],\n  [\"965\u00d77=\", \"333\u00d79=\"],\n  [\"407\u00d74=\", \"771\u00d75=\"],\n  [\"195\u00d78=\", \"466\u00d79=\"],\n  [\"186\u00d74=\", \"551\u00d73=\"],\n  [\"807\u00d77=\", \"503\u00d74=\"],\n  [\"479\u00d79=\", \"282\u00d76=\"],\n  [\"566\u00d78=\", \"863\u00d77=\"],\n  [\"857\u00d74=\", \"371\u00d75=\"],\n  [\"716\u00d78=\", \"684\u00d73=\"],\n  [\"794\u00d75=\", \"690\u00d75=\"],\n  [\"906\u00d78=\", \"841\u00d76=\"],\n  [\"887\u00d75=\", \"482\u00d79=\"],\n  [\"616\u00d72=\", \"170\u00d72=\"],\n  [\"868\u00d79=\", \"192\u00d76=\"],\n  [\"944\u00d75=\", \"524\u00d72=\"],\n  [\"632\u00d79=\", \"119\u00d79=\"],\n  [\"133\u00d77=\", \"507\u00d74=\"],\n  [\"859\u00d72=\", \"316\u00d74=\"],\n  [\"407\u00d75=\", \"316\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 24 multiplication expressions in the table,\n# per the authoring diff. Every \"old\" string is unique within the document,\n# so Find/Execute locates exactly one occurrence each; wdReplaceAll (2) is\n# used defensively but only ever touches the single match. Run formatting\n# (fonts/size) is preserved because Find/Replace edits text in place.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-10-10 Friday\", \"2025-10-11 Saturday\"),\n  @(\"887\u00d78=\", \"201\u00d76=\"),\n  @(\"953\u00d75=\", \"257\u00d76=\"),\n  @(\"463\u00d74=\", \"411\u00d77=\"),\n  @(\"141\u00d74=\", \"191\u00d72=\"),\n  @(\"579\u00d75=\", \"580\u00d72=\"),\n  @(\"169\u00d75=\", \"643\u00d76=\"),\n  @(\"965\u00d77=\", \"333\u00d79=\"),\n  @(\"407\u00d74=\", \"771\u00d75=\"),\n  @(\"195\u00d78=\", \"466\u00d79=\"),\n  @(\"186\u00d74=\", \"551\u00d73=\"),\n  @(\"807\u00d77=\", \"503\u00d74=\"),\n  @(\"479\u00d79=\", \"282\u00d76=\"),\n  @(\"566\u00d78=\", \"863\u00d77=\"),\n  @(\"857\u00d74=\", \"371\u00d75=\"),\n  @(\"716\u00d78=\", \"684\u00d73=\"),\n  @(\"794\u00d75=\", \"690\u00d75=\"),\n  @(\"906\u00d78=\", \"841\u00d76=\"),\n  @(\"887\u00d75=\", \"482\u00d79=\"),\n  @(\"616\u00d72=\", \"170\u00d72=\"),\n  @(\"868\u00d79=\", \"192\u00d76=\"),\n  @(\"944\u00d75=\", \"524\u00d72=\"),\n  @(\"632\u00d79=\", \"119\u00d79=\"),\n  @(\"133\u00d77=\", \"507\u00d74=\"),\n  @(\"859\u00d72=\", \"316\u00d74=\"),\n  @(\"407\u00d75=\", \"316\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
